$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B14").Value = "MissingNumber"
$ws.Range("A14").Value = "Find the Missing Number"

$ws.Range("A14").Select()
